# Ran code for averaged intensities on spiral schemes
# Rows 10-16 (existing "HKL" indices 8-14) are refreshed with newly computed
# averaged-intensity values, and the scheme previously associated with each
# row slides down to make room for three brand-new spiral sampling schemes
# ("Spiral-90deg-10rot-5space", "Spiral-90deg-15rot-5space",
# "Spiral-90deg-10rot-3space"). Three new rows (17-19) are appended holding
# the data that used to live in rows 13-15 (HexGrid-* schemes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Gaussian-Quadrature -------------------------------------------
$ws.Cells.Item(10,2).Value = "Gaussian-Quadrature"
$ws.Cells.Item(10,3).Value = 1.201031381252061
$ws.Cells.Item(10,4).Value = 0.9924703951312377
$ws.Cells.Item(10,5).Value = 0.9834174039073463
$ws.Cells.Item(10,6).Value = 0.9305288542310139
$ws.Cells.Item(10,7).Value = 1.201031381252061
$ws.Cells.Item(10,8).Value = 0.9924703951312377
$ws.Cells.Item(10,9).Value = 1.036209067011093
$ws.Cells.Item(10,10).Value = 0.9172375850247387
$ws.Cells.Item(10,11).Value = 1.026900555561042
$ws.Cells.Item(10,12).Value = 0.9346251806368507
$ws.Cells.Item(10,13).Value = 1.201031381252061
$ws.Cells.Item(10,14).Value = 0.987943899519292
$ws.Cells.Item(10,15).Value = 1.026862008630415
$ws.Cells.Item(10,16).Value = 1.002802552844423

# --- Row 11: Spiral-90deg-10rot-5space -------------------------------------
$ws.Cells.Item(11,2).Value = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(11,3).Value = 0.328386906376964
$ws.Cells.Item(11,4).Value = 1.889623144686522
$ws.Cells.Item(11,5).Value = 0.9401723248035061
$ws.Cells.Item(11,6).Value = 1.112288761701969
$ws.Cells.Item(11,7).Value = 0.328386906376964
$ws.Cells.Item(11,8).Value = 1.889623144686522
$ws.Cells.Item(11,9).Value = 0.725773798042578
$ws.Cells.Item(11,10).Value = 1.206432846945116
$ws.Cells.Item(11,11).Value = 0.7018675165841424
$ws.Cells.Item(11,12).Value = 1.535077249173352
$ws.Cells.Item(11,13).Value = 0.328386906376964
$ws.Cells.Item(11,14).Value = 1.414897734745014
$ws.Cells.Item(11,15).Value = 1.06761778439224
$ws.Cells.Item(11,16).Value = 1.054952818539269

# --- Row 12: Spiral-90deg-15rot-5space -------------------------------------
$ws.Cells.Item(12,2).Value = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(12,3).Value = 0.3274990554223933
$ws.Cells.Item(12,4).Value = 1.893761146967791
$ws.Cells.Item(12,5).Value = 0.9374496479028154
$ws.Cells.Item(12,6).Value = 1.113453893430607
$ws.Cells.Item(12,7).Value = 0.3274990554223933
$ws.Cells.Item(12,8).Value = 1.893761146967791
$ws.Cells.Item(12,9).Value = 0.7235663637782129
$ws.Cells.Item(12,10).Value = 1.206605241917441
$ws.Cells.Item(12,11).Value = 0.7015036444599082
$ws.Cells.Item(12,12).Value = 1.537951313656142
$ws.Cells.Item(12,13).Value = 0.3274990554223933
$ws.Cells.Item(12,14).Value = 1.415605397435303
$ws.Cells.Item(12,15).Value = 1.068040935930902
$ws.Cells.Item(12,16).Value = 1.055223788441914

# --- Row 13: Spiral-90deg-10rot-3space -------------------------------------
$ws.Cells.Item(13,2).Value = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(13,3).Value = 0.3282244278643351
$ws.Cells.Item(13,4).Value = 1.890164067852492
$ws.Cells.Item(13,5).Value = 0.9396232570599472
$ws.Cells.Item(13,6).Value = 1.112596242688806
$ws.Cells.Item(13,7).Value = 0.3282244278643351
$ws.Cells.Item(13,8).Value = 1.890164067852492
$ws.Cells.Item(13,9).Value = 0.7252233033950867
$ws.Cells.Item(13,10).Value = 1.206515316486272
$ws.Cells.Item(13,11).Value = 0.7018487639650598
$ws.Cells.Item(13,12).Value = 1.535686761862671
$ws.Cells.Item(13,13).Value = 0.3282244278643351
$ws.Cells.Item(13,14).Value = 1.414893662456219
$ws.Cells.Item(13,15).Value = 1.067651998866395
$ws.Cells.Item(13,16).Value = 1.054985267646833

# --- Row 14: NoRotation-tilt60deg (was row 10's scheme) ---------------------
$ws.Cells.Item(14,2).Value = "NoRotation-tilt60deg"
$ws.Cells.Item(14,3).Value = 0.4694839999999998
$ws.Cells.Item(14,4).Value = 1.190688000000001
$ws.Cells.Item(14,5).Value = 1.408632
$ws.Cells.Item(14,6).Value = 0.9236320000000016
$ws.Cells.Item(14,7).Value = 0.4694839999999998
$ws.Cells.Item(14,8).Value = 1.190688000000001
$ws.Cells.Item(14,9).Value = 1.059008000000001
$ws.Cells.Item(14,10).Value = 1.177095999999999
$ws.Cells.Item(14,11).Value = 0.7513119999999998
$ws.Cells.Item(14,12).Value = 1.101719999999999
$ws.Cells.Item(14,13).Value = 0.4694839999999998
$ws.Cells.Item(14,14).Value = 1.299660000000001
$ws.Cells.Item(14,15).Value = 0.9981090000000006
$ws.Cells.Item(14,16).Value = 1.0101965

# --- Row 15: Rotation-NoTilt (was row 11's scheme) --------------------------
$ws.Cells.Item(15,2).Value = "Rotation-NoTilt"
$ws.Cells.Item(15,3).Value = 0.66
$ws.Cells.Item(15,4).Value = 0.21
$ws.Cells.Item(15,5).Value = 1.97
$ws.Cells.Item(15,6).Value = 0.7
$ws.Cells.Item(15,7).Value = 0.66
$ws.Cells.Item(15,8).Value = 0.21
$ws.Cells.Item(15,9).Value = 1.498749999999999
$ws.Cells.Item(15,10).Value = 1.160275000000001
$ws.Cells.Item(15,11).Value = 0.8426125000000002
$ws.Cells.Item(15,12).Value = 0.4952749999999999
$ws.Cells.Item(15,13).Value = 0.66
$ws.Cells.Item(15,14).Value = 1.09
$ws.Cells.Item(15,15).Value = 0.885
$ws.Cells.Item(15,16).Value = 0.9421140625

# --- Row 16: Rotation-60detTilt (was row 12's scheme) -----------------------
$ws.Cells.Item(16,2).Value = "Rotation-60detTilt"
$ws.Cells.Item(16,3).Value = 0.8104218732544011
$ws.Cells.Item(16,4).Value = 0.5397915715584015
$ws.Cells.Item(16,5).Value = 1.5545899380736
$ws.Cells.Item(16,6).Value = 0.8263806980096018
$ws.Cells.Item(16,7).Value = 0.8104218732544011
$ws.Cells.Item(16,8).Value = 0.5397915715584015
$ws.Cells.Item(16,9).Value = 1.281718131097599
$ws.Cells.Item(16,10).Value = 1.0857705370624
$ws.Cells.Item(16,11).Value = 0.9099897923583969
$ws.Cells.Item(16,12).Value = 0.7063428915200011
$ws.Cells.Item(16,13).Value = 0.8104171606016012
$ws.Cells.Item(16,14).Value = 1.047190754816001
$ws.Cells.Item(16,15).Value = 0.9327960202240011
$ws.Cells.Item(16,16).Value = 0.9643756791168002

# --- New row 17: HexGrid-90degTilt5degRes (was row 13's scheme) -------------
$ws.Cells.Item(16,1).Copy()
$ws.Cells.Item(17,1).PasteSpecial(-4122)
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(17,3).Value = 0.9936876058048703
$ws.Cells.Item(17,4).Value = 0.9933137835985731
$ws.Cells.Item(17,5).Value = 0.9954199582472592
$ws.Cells.Item(17,6).Value = 0.9938893334660253
$ws.Cells.Item(17,7).Value = 0.9936876058048703
$ws.Cells.Item(17,8).Value = 0.9933137835985731
$ws.Cells.Item(17,9).Value = 0.9943988124178552
$ws.Cells.Item(17,10).Value = 0.9961554661360671
$ws.Cells.Item(17,11).Value = 0.9942538958394512
$ws.Cells.Item(17,12).Value = 0.9920335538727214
$ws.Cells.Item(17,13).Value = 0.9936943093737369
$ws.Cells.Item(17,14).Value = 0.9943668709229161
$ws.Cells.Item(17,15).Value = 0.9940776702791819
$ws.Cells.Item(17,16).Value = 0.9941440511728529

# --- New row 18: HexGrid-90degTilt22p5degRes (was row 14's scheme) ----------
$ws.Cells.Item(16,1).Copy()
$ws.Cells.Item(18,1).PasteSpecial(-4122)
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = "HexGrid-90degTilt22p5degRes"
$ws.Cells.Item(18,3).Value = 1.002644322090684
$ws.Cells.Item(18,4).Value = 1.02901070469991
$ws.Cells.Item(18,5).Value = 0.9814080470948402
$ws.Cells.Item(18,6).Value = 0.9915027306285648
$ws.Cells.Item(18,7).Value = 1.002644322090684
$ws.Cells.Item(18,8).Value = 1.02901070469991
$ws.Cells.Item(18,9).Value = 0.9835592198750354
$ws.Cells.Item(18,10).Value = 0.991941573237313
$ws.Cells.Item(18,11).Value = 0.9904784264816872
$ws.Cells.Item(18,12).Value = 1.00702966372145
$ws.Cells.Item(18,13).Value = 1.002644322090684
$ws.Cells.Item(18,14).Value = 1.005209375897375
$ws.Cells.Item(18,15).Value = 1.0011414511285
$ws.Cells.Item(18,16).Value = 0.9971968359786855

# --- New row 19: HexGrid-60degTilt5degRes (was row 15's scheme) ------------
$ws.Cells.Item(16,1).Copy()
$ws.Cells.Item(19,1).PasteSpecial(-4122)
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(19,3).Value = 0.9738047511993261
$ws.Cells.Item(19,4).Value = 1.069272647050894
$ws.Cells.Item(19,5).Value = 0.976643070494326
$ws.Cells.Item(19,6).Value = 1.004721176537816
$ws.Cells.Item(19,7).Value = 0.9738047511993261
$ws.Cells.Item(19,8).Value = 1.069272647050894
$ws.Cells.Item(19,9).Value = 0.9685981757026126
$ws.Cells.Item(19,10).Value = 0.9969514031025885
$ws.Cells.Item(19,11).Value = 0.9815719998120205
$ws.Cells.Item(19,12).Value = 1.036769794439071
$ws.Cells.Item(19,13).Value = 0.9738182289009422
$ws.Cells.Item(19,14).Value = 1.02295785877261
$ws.Cells.Item(19,15).Value = 1.00611041132059
$ws.Cells.Item(19,16).Value = 1.001041627292332
